# Fixed calculate fitness functions.
# Remove a handful of instructors, add a new "Çalışma Günleri" value for every
# remaining instructor, and shrink the data table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final list of remaining instructors, in their new (alphabetically sorted) order,
# each working Monday through Friday.
$workDays = "Pazartesi, Salı, Çarşamba, Perşembe, Cuma"

$names = @(
    "Aydın Seçer",
    "Didem Yılmaz Çapkur",
    "Elham Pashaei",
    "Fatih Koçan",
    "Gülsüm Yeliz Şentürk",
    "Hakan Aydın",
    "Ümit Alkan",
    "Kenan Özden",
    "Tarık Çakar"
)

# Rewrite rows 2..10 (ID, İsim, İş Günleri) with the final data set.
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $workDays
}

# Clear out the old rows 11..15 that are no longer part of the data set.
$ws.Range("A11:C15").ClearContents()

# Shrink the table to match the new data extent.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C10"))

# Update the selected cell as in the saved workbook.
$ws.Range("J7").Select()
